# Update TPM-derived NATMI metrics on the active worksheet.
# Columns M..T hold Receptor/Edge expression values and specificities that
# were recomputed using new TPM values; only the affected cells are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 0.1318220377821111
$ws.Range("R2").Value = 1.186398340039
$ws.Range("S2").Value = 0.002607375069995422
$ws.Range("T2").Value = 0.002607375069995421

# Row 3
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("S3").Value = 0.6476063272730862
$ws.Range("T3").Value = 0.6476063272730862

# Row 4
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("Q4").Value = 4.794432716175666
$ws.Range("S4").Value = 0.09483152096002037
$ws.Range("T4").Value = 0.09483152096002037

# Row 5
$ws.Range("M5").Value = 0.4102596666666667
$ws.Range("N5").Value = 1.230779
$ws.Range("O5").Value = 0.003499619873322347
$ws.Range("P5").Value = 0.003499619873322347
$ws.Range("Q5").Value = 0.04510955463544445
$ws.Range("R5").Value = 0.405985991719
$ws.Range("S5").Value = 0.0008922448033269254
$ws.Range("T5").Value = 0.0008922448033269253

# Row 6
$ws.Range("O6").Value = 0.8692174743460166
$ws.Range("P6").Value = 0.8692174743460165
$ws.Range("S6").Value = 0.2216111470729303
$ws.Range("T6").Value = 0.2216111470729303

# Row 7
$ws.Range("N7").Value = 44.764041
$ws.Range("O7").Value = 0.1272829057806611
$ws.Range("P7").Value = 0.1272829057806611
$ws.Range("S7").Value = 0.03245138482064077
$ws.Range("T7").Value = 0.03245138482064077
